# Mas datos del CFD del DBX 1. Hay controles elevator hasta 25 deg AoA
# Adds three new elevator-sweep CFD cases (AoA 22, 25, 35 deg) in the
# symmetric ("SI") block at the bottom of the table, and underlines the
# whole "Simetria" (column B) data range to flag it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace old row 19 (dbx..a30_bX..) with the new elevator sweep cases,
# and populate the previously-empty row 21.
$ws.Range("A19").Value = "dbx_v1_s50_r16_a22_b0_da0_deX_dr0"
$ws.Range("A20").Value = "dbx_v1_s50_r16_a25_b0_da0_deX_dr0"
$ws.Range("A21").Value = "dbx_v1_s50_r16_a35_b0_da0_deX_dr0"

$ws.Range("B19").Value = "SI"
$ws.Range("B20").Value = "SI"
$ws.Range("B21").Value = "SI"

# Row 21 previously had no formatting (blank cells) - copy the formatting
# from the row above so it matches the rest of the table (border/shading).
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)  # xlPasteFormats

# Underline the whole "Simetria" column data range to highlight it.
$ws.Range("B2:B21").Font.Underline = $true

# Update the active selection to the underlined range.
$ws.Range("B2").Select()
$ws.Range("B2:B21").Select()
